# updated main GSC export data
# Appends the newest day's row (2025-12-25) to the "Chart" sheet, which
# drives the daily export data (Date / Invalid / Valid columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Insert a new row right after the current last data row (row 81) so the
# sheet grows from A1:C81 to A1:C82, same as adding one more day's record.
$ws.Rows.Item(82).Insert()

# Column A holds dates as literal text (e.g. "2025-12-24"), not real date
# serials -- force a text number format before assigning the string so
# Excel doesn't auto-coerce it into a date value, then drop the format
# again so the cell matches the plain (unformatted) style used by every
# other row in the column.
$ws.Range("A82").NumberFormat = "@"
$ws.Range("A82").Value = "2025-12-25"
$ws.Range("A82").ClearFormats()

# Invalid / Valid counts for the new day.
$ws.Range("B82").Value = 0
$ws.Range("C82").Value = 32
